$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected; unprotect to apply edits, then re-protect at the end.
$ws.Unprotect()

# Update the confidential disclaimer date string (A18): 2021-03-18 -> 2021-03-19
$ws.Range("A18").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-19 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values for rows 2-15
$ws.Range("D2").Value = 0.05580303976851266
$ws.Range("E2").Value = 0.00126829268292683

$ws.Range("D3").Value = 0.02361353264545879
$ws.Range("E3").Value = 0.003429493645350057

$ws.Range("D4").Value = 0.03202618115299968
$ws.Range("E4").Value = 0.00994454006502199

$ws.Range("D5").Value = 0.03214391195641374
$ws.Range("E5").Value = 0.001206757843926098

$ws.Range("D6").Value = 0.034200458129346
$ws.Range("E6").Value = -0.004536771728749001

$ws.Range("D7").Value = 0.01923095244900487
$ws.Range("E7").Value = -0.004458756502353167

$ws.Range("D8").Value = 0.004642540698790163
$ws.Range("E8").Value = 0.02638522427440626

$ws.Range("D9").Value = 0.006610618637943072
$ws.Range("E9").Value = -0.002676549310273835

$ws.Range("D10").Value = 0.07043432458586658
$ws.Range("E10").Value = 0.003478260869565153

$ws.Range("D11").Value = 0.07051598757089368
$ws.Range("E11").Value = 0.004632310364794501

$ws.Range("D12").Value = 0.1458174260643749
$ws.Range("E12").Value = 0.00619772998805268

$ws.Range("D13").Value = 0.3925335532789729
$ws.Range("E13").Value = -0.0006189213085764367

$ws.Range("D14").Value = 0.1124274730614228
$ws.Range("E14").Value = 0.01144019321215195

$ws.Range("D15").Value = 0.9999999999999998
$ws.Range("E15").Value = 0.002891550194834114

# Restore sheet protection
$ws.Protect()
